$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (which already has the bold/bordered/centered
# header style) onto the two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(7, 8),
    @(5, 6),
    @(4, 5),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(5, 6),
    @(6, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(5, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
